$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = "UNFI – Dublin, CA Appt:"
$ws.Range("E15").Value = "dubappts@unfi.com"
$ws.Range("E16").Value = "4000 Inspiration Drive"
$ws.Range("E17").Value = "Dublin, CA 94568"
$ws.Range("E21").Value = ""
$ws.Range("C26").Value = "10827"
$ws.Range("C27").Value = "10829"
$ws.Range("C28").Value = "82058"
$ws.Range("C29").Value = "11311"
